$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 2.02
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 3.4
$ws.Range("L2").Value = 1.51
$ws.Range("N2").Value = 2.4
$ws.Range("O2").Value = 1.58
$ws.Range("P2").Value = 1.46
$ws.Range("Q2").Value = 2.74
$ws.Range("S2").Value = 5.9
$ws.Range("T2").Value = 2.26
$ws.Range("V2").Value = 1.24
$ws.Range("W2").Value = 1.82
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 990
$ws.Range("AB2").Value = 12
$ws.Range("AC2").Value = 7.8
$ws.Range("AF2").Value = 32
$ws.Range("AG2").Value = 28
$ws.Range("AK2").Value = 140
$ws.Range("AL2").Value = 210
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 3.45
$ws.Range("I3").Value = 3.05
$ws.Range("O3").Value = 1.65
$ws.Range("R3").Value = 1.14
$ws.Range("S3").Value = 6.4
$ws.Range("T3").Value = 2.26
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 1.41
$ws.Range("AB3").Value = 22
$ws.Range("AC3").Value = 22
$ws.Range("F4").Value = 1.16
$ws.Range("G4").Value = 1.17
$ws.Range("I4").Value = 27
$ws.Range("K4").Value = 11
$ws.Range("N4").Value = 9.4
$ws.Range("O4").Value = 1.1
$ws.Range("P4").Value = 3.8
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.72
$ws.Range("T4").Value = 1.96
$ws.Range("U4").Value = 1.86
$ws.Range("W4").Value = 6.8
$ws.Range("X4").Value = 1000
$ws.Range("Z4").Value = 290
$ws.Range("AB4").Value = 16
$ws.Range("AC4").Value = 40
$ws.Range("AD4").Value = 990
$ws.Range("AE4").Value = 380
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 990
$ws.Range("AJ4").Value = 9.800000000000001
$ws.Range("AK4").Value = 14
$ws.Range("AL4").Value = 36
$ws.Range("G6").Value = 2.62
$ws.Range("S6").Value = 4.8
$ws.Range("H7").Value = 5.3
$ws.Range("I7").Value = 5.9
$ws.Range("U7").Value = 2.68
$ws.Range("V7").Value = 1.2
$ws.Range("Y7").Value = 990
$ws.Range("AC7").Value = 14.5
$ws.Range("AD7").Value = 24
$ws.Range("AF7").Value = 15.5
$ws.Range("AH7").Value = 18
$ws.Range("AI7").Value = 110
$ws.Range("AJ7").Value = 18
$ws.Range("AM7").Value = 55
$ws.Range("AN7").Value = 5.1
$ws.Range("G8").Value = 5.9
$ws.Range("H8").Value = 1.74
$ws.Range("V8").Value = 2.18
$ws.Range("AA8").Value = 38
$ws.Range("AB8").Value = 1000
$ws.Range("AE8").Value = 65
$ws.Range("AI8").Value = 95
$ws.Range("F9").Value = 1.63
$ws.Range("G9").Value = 1.73
$ws.Range("W9").Value = 2.36
$ws.Range("F10").Value = 1.44
$ws.Range("G10").Value = 1.45
$ws.Range("N10").Value = 3.55
$ws.Range("O10").Value = 1.37
$ws.Range("P10").Value = 1.89
$ws.Range("R10").Value = 1.33
$ws.Range("T10").Value = 2.38
$ws.Range("W10").Value = 3.2
$ws.Range("X10").Value = 14
$ws.Range("M11").Value = 1.05
$ws.Range("T11").Value = 1.87
$ws.Range("X11").Value = 18.5
$ws.Range("Y11").Value = 23
$ws.Range("AH11").Value = 20
$ws.Range("AI11").Value = 75
$ws.Range("AO11").Value = 100
$ws.Range("J12").Value = 3.1
$ws.Range("N12").Value = 2.88
$ws.Range("U12").Value = 1.87
$ws.Range("AM12").Value = 180